$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.072.92'
$ws.Range("E2").Value = '  +1.34%  '
$ws.Range("D3").Value = '3.171.52'
$ws.Range("E3").Value = '  +3.47%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.79%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.170.23'
$ws.Range("E8").Value = '  +3.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.163'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.22'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.72%  '
$ws.Range("E12").Value = '  +4.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000274'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +17.89%  '
$ws.Range("E14").Value = '  +6.31%  '
$ws.Range("D15").Value = '3.687.30'
$ws.Range("E15").Value = '  +3.51%  '
$ws.Range("D16").Value = '65.148.14'
$ws.Range("E16").Value = '  +1.51%  '
$ws.Range("D17").Value = '3.167.99'
$ws.Range("E17").Value = '  +3.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.03%  '
$ws.Range("E19").Value = '  +1.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '511.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.40%  '
$ws.Range("B22").Value = 'Polygon'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.734'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.47%  '
$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.65'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.85'
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.52%  '
$ws.Range("E28").Value = '  +3.35%  '
$ws.Range("E29").Value = '  +7.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '28.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.23%  '
$ws.Range("E31").Value = '  +14.39%  '
$ws.Range("E32").Value = '  +7.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.33'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +10.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.68'
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = '  +1.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '481.47'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0887'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.09'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.44%  '
$ws.Range("E40").Value = '  +2.59%  '
$ws.Range("D41").Value = '3.125.30'
$ws.Range("E41").Value = '  +4.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.63'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.24%  '
$ws.Range("E43").Value = '  +4.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +16.55%  '
$ws.Range("E45").Value = '  +10.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.15%  '
$ws.Range("D47").Value = '0.0₃0597'
$ws.Range("E47").Value = '  +14.43%  '
$ws.Range("E49").Value = '  +1.51%  '
$ws.Range("E50").Value = '  +10.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '122.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.00%  '
